$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text format so numeric-looking strings (e.g. "576.43")
# are not silently converted to numbers, losing exact text representation
# (trailing zeros, multi-dot thousand separators, etc.)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "64.995.24"
$ws.Range("E2").Value = "  +1.67%  "
$ws.Range("D3").Value = "3.154.91"
$ws.Range("E3").Value = "  +2.98%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "576.43"
$ws.Range("E5").Value = "  +3.11%  "
$ws.Range("D6").Value = "149.78"
$ws.Range("E6").Value = "  +4.85%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "3.152.21"
$ws.Range("E8").Value = "  +3.01%  "
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("E10").Value = "  +4.10%  "
$ws.Range("D11").Value = "6.12"
$ws.Range("E11").Value = "  -1.37%  "
$ws.Range("D12").Value = "0.500"
$ws.Range("E12").Value = "  +3.92%  "
$ws.Range("D13").Value = "0.0000264"
$ws.Range("E13").Value = "  +13.55%  "
$ws.Range("D14").Value = "37.12"
$ws.Range("E14").Value = "  +5.16%  "
$ws.Range("D15").Value = "3.676.08"
$ws.Range("E15").Value = "  +2.99%  "
$ws.Range("D16").Value = "65.094.38"
$ws.Range("E16").Value = "  +1.74%  "
$ws.Range("D17").Value = "3.177.03"
$ws.Range("D18").Value = "7.10"
$ws.Range("E18").Value = "  +4.69%  "
$ws.Range("E19").Value = "  +1.22%  "
$ws.Range("D20").Value = "509.32"
$ws.Range("E20").Value = "  +4.65%  "
$ws.Range("D21").Value = "14.80"
$ws.Range("E21").Value = "  +3.43%  "
$ws.Range("D23").Value = "15.26"
$ws.Range("E23").Value = "  +4.39%  "
$ws.Range("D24").Value = "7.73"
$ws.Range("E24").Value = "  +2.38%  "
$ws.Range("D25").Value = "84.41"
$ws.Range("E25").Value = "  +2.03%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("E27").Value = "  +3.44%  "
$ws.Range("E28").Value = "  +8.97%  "
$ws.Range("D29").Value = "2.16"
$ws.Range("E29").Value = "  +5.54%  "
$ws.Range("D30").Value = "27.64"
$ws.Range("E30").Value = "  +4.05%  "
$ws.Range("D31").Value = "2.77"
$ws.Range("E31").Value = "  +9.85%  "
$ws.Range("D32").Value = "0.999"
$ws.Range("E32").Value = "  -0.09%  "
$ws.Range("E33").Value = "  +2.36%  "
$ws.Range("D34").Value = "6.23"
$ws.Range("E34").Value = "  +9.36%  "
$ws.Range("D35").Value = "6.53"
$ws.Range("E35").Value = "  +4.70%  "
$ws.Range("D36").Value = "55.21"
$ws.Range("E36").Value = "  +0.40%  "
$ws.Range("E37").Value = "  +9.87%  "
$ws.Range("E38").Value = "  +3.10%  "
$ws.Range("D39").Value = "464.46"
$ws.Range("E39").Value = "  +4.40%  "
$ws.Range("D40").Value = "3.00"
$ws.Range("E40").Value = "  +6.87%  "
$ws.Range("D41").Value = "8.66"
$ws.Range("E41").Value = "  +3.92%  "
$ws.Range("D42").Value = "3.067.37"
$ws.Range("E42").Value = "  +1.35%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "0.282"
$ws.Range("E44").Value = "  +1.91%  "
$ws.Range("E45").Value = "  +6.78%  "
$ws.Range("D46").Value = "28.65"
$ws.Range("E46").Value = "  +3.11%  "
$ws.Range("E47").Value = "  +13.25%  "
$ws.Range("E48").Value = "  +0.04%  "
$ws.Range("E49").Value = "  +0.50%  "
$ws.Range("E50").Value = "  +5.46%  "
$ws.Range("D51").Value = "119.92"
$ws.Range("E51").Value = "  +1.71%  "

# Restore the original (default) cell formatting now that the text values
# are safely stored, so no stray style attributes remain on the cells.
$ws.Range("D2:D51").ClearFormats()
